# Add function to include member name in measurement report.
# This adds two new measurement rows (7 & 8) to the "Main" sheet, the second
# of which is the first entry to carry a "Team member" (column J) value -
# i.e. the new "include member name" data point - and nudges the two
# timestamp cells on row 6 that got re-stamped by the same recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: timestamps re-saved with slightly different sub-millisecond
#     precision (same instant, just re-serialised) ---
$ws.Range("C6").Value = 45406.94200855324
$ws.Range("D6").Value = 45406.94212322916

# --- Row 7: new measurement, no team member recorded yet ---
$ws.Range("A7").Value = 6
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "Snakes"
$ws.Range("B7").Style = "Normal"

$ws.Range("C7").Value = 45406.94871221065
$ws.Range("C7").Style = "Normal"
$ws.Range("C7").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("D7").Value = 45406.94883653936
$ws.Range("D7").Style = "Normal"
$ws.Range("D7").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("F7").Value = 0.0001157407407407407
$ws.Range("F7").Style = "Normal"
$ws.Range("F7").NumberFormat = "h:mm:ss"

$ws.Range("H7").Value = "Team3"
$ws.Range("H7").Style = "Normal"

$ws.Range("I7").Value = "Process13"
$ws.Range("I7").Style = "Normal"

$ws.Range("J7").Value = "Process13"
$ws.Range("J7").Style = "Normal"

# --- Row 8: new measurement, now including the team member name ---
$ws.Range("A8").Value = 7
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = "Snakes"
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = 45406.95036223251
$ws.Range("C8").Style = "Normal"
$ws.Range("C8").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("D8").Value = 45406.95041998194
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("F8").Value = 0.00004629629629629629
$ws.Range("F8").Style = "Normal"
$ws.Range("F8").NumberFormat = "h:mm:ss"

$ws.Range("H8").Value = "Team1"
$ws.Range("H8").Style = "Normal"

$ws.Range("I8").Value = "Process1"
$ws.Range("I8").Style = "Normal"

$ws.Range("J8").Value = "Person1"
$ws.Range("J8").Style = "Normal"
